# Update X_CHR values in the "Individual Costs" sheet with new prices.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individual Costs")

$newValues = @(
    242.9301789682836,
    241.7533522765146,
    244.4072642249967,
    230.6829339348994,
    219.1516994869718,
    211.218754640267,
    203.3390079515765,
    201.3797721673737,
    245.4722048135623,
    254.8313851478964,
    260.2003103654174,
    259.9537246242008,
    252.9227432672594,
    267.3152168258129,
    300.3356876444313,
    305.7807857621891,
    258.7384440859091,
    237.7027714931847,
    281.3168672212262,
    256.801896043078,
    238.1904045372331,
    229.3050071278335,
    203.1075069825606,
    217.5450755808948
)

$row = 2
foreach ($val in $newValues) {
    $ws.Cells.Item($row, 3).Value = $val
    $row++
}
